# Generate Report for Handoff
# - Remove the row for 54f9528b-44ec-45a5-9ec9-d8b507c5d970.md (now resolved / out of report)
# - Promote the row for f7c573b5-210d-4457-aa12-dc21f8919674.md to "Ready for handoff"
#   with refreshed handoff/handback timestamps and an out-of-date warning message.
# - Widen the "Error Detail" column on the per-locale sheets to fit the new message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop the hyperlinks before the row shuffle - they get rebuilt from scratch below.
$ov.Range("B2").Hyperlinks.Delete()

# Remove the 54f9528b row entirely; f7c573b5's row shifts up from row 4 to row 3.
$ov.Range("A3:G3").EntireRow.Delete()

# Refresh the (now row 3) f7c573b5 status columns.
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-11-10 06:42:00"

$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0ea60244a7b1ec7d192c57cd9e5cb447e8e75ef/e2e/4b740145-cd67-40db-9228-17010fb6e6a3.md", "", "", "e2e\4b740145-cd67-40db-9228-17010fb6e6a3.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0ea60244a7b1ec7d192c57cd9e5cb447e8e75ef/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md", "", "", "e2e\f7c573b5-210d-4457-aa12-dc21f8919674.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Hyperlinks.Delete()

$zh.Range("A3:P3").EntireRow.Delete()

$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("H3").Value = "2016-11-10 06:41:45"
$zh.Range("K3").Value = "2016-11-10 06:38:45"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0ea60244a7b1ec7d192c57cd9e5cb447e8e75ef/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3597f57f0d8c423d1fff5bf482dd93ca333122bd/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md."

$zh.Columns.Item(16).ColumnWidth = 40

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0ea60244a7b1ec7d192c57cd9e5cb447e8e75ef/e2e/4b740145-cd67-40db-9228-17010fb6e6a3.md", "", "", "4b740145-cd67-40db-9228-17010fb6e6a3.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b793f56d6982ece26bdea7bb7a73f49f65c528cd/e2e/4b740145-cd67-40db-9228-17010fb6e6a3.md", "", "", "4b740145-cd67-40db-9228-17010fb6e6a3.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0ea60244a7b1ec7d192c57cd9e5cb447e8e75ef/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md", "", "", "f7c573b5-210d-4457-aa12-dc21f8919674.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b793f56d6982ece26bdea7bb7a73f49f65c528cd/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md", "", "", "f7c573b5-210d-4457-aa12-dc21f8919674.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Hyperlinks.Delete()

$de.Range("A3:P3").EntireRow.Delete()

$de.Range("C3").Value = "Ready for handoff"
$de.Range("H3").Value = "2016-11-10 06:42:00"
$de.Range("K3").Value = "2016-11-10 06:39:05"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0ea60244a7b1ec7d192c57cd9e5cb447e8e75ef/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3597f57f0d8c423d1fff5bf482dd93ca333122bd/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md."

$de.Columns.Item(16).ColumnWidth = 40

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0ea60244a7b1ec7d192c57cd9e5cb447e8e75ef/e2e/4b740145-cd67-40db-9228-17010fb6e6a3.md", "", "", "4b740145-cd67-40db-9228-17010fb6e6a3.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2371aaff29d75d63dec0b03ef73126cee1c5f0e4/e2e/4b740145-cd67-40db-9228-17010fb6e6a3.md", "", "", "4b740145-cd67-40db-9228-17010fb6e6a3.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0ea60244a7b1ec7d192c57cd9e5cb447e8e75ef/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md", "", "", "f7c573b5-210d-4457-aa12-dc21f8919674.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2371aaff29d75d63dec0b03ef73126cee1c5f0e4/e2e/f7c573b5-210d-4457-aa12-dc21f8919674.md", "", "", "f7c573b5-210d-4457-aa12-dc21f8919674.md") | Out-Null
